$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Row 7 ---------------------------------------------------------------
$ws.Range("B7").Value = 2386
$ws.Range("C7").Value = 6365
$ws.Range("D7").Formula = "=3*10^8/(C7*10^(-10))"
$ws.Range("E7").Formula = "=D7*6.62*10^(-34)/(1.6*10^(-19))"
$ws.Range("F7").Formula = "=E7*0.01"
$ws.Range("I7").Formula = "=E8-E7"
$ws.Range("J7").Formula = "=F8+F7"

# --- Row 8 -----------------------------------------------------------------
$ws.Range("B8").Value = 2282
$ws.Range("C8").Value = 6110
$ws.Range("D8").Formula = "=3*10^8/(C8*10^(-10))"
$ws.Range("E8").Formula = "=D8*6.62*10^(-34)/(1.6*10^(-19))"
$ws.Range("F8").Formula = "=E8*0.01"
$ws.Range("I8").Formula = "=I7/5"
$ws.Range("J8").Formula = "=J7/5"

# --- Row 9 -------------------------------------------------------------
$ws.Range("B9").Value = 1616
$ws.Range("C9").Value = 5023
$ws.Range("D9").Formula = "=3*10^8/(C9*10^(-10))"
$ws.Range("E9").Formula = "=D9*6.62*10^(-34)/(1.6*10^(-19))"
$ws.Range("F9").Formula = "=E9*0.01"

# --- Formatting: copy the cell style used by the original data rows ------
# (style index 1 - "Helvetica Neue" 10pt) onto the new numeric cells that
# should carry it, leaving C9 with the default style (matches source data).
$ws.Range("B2:B4").Copy() | Out-Null
$ws.Range("B7:B9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C2:C3").Copy() | Out-Null
$ws.Range("C7:C8").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Re-apply the values (PasteSpecial(xlPasteFormats) does not touch values,
# but keep this to guarantee correctness regardless of paste semantics).
$ws.Range("B7").Value = 2386
$ws.Range("B8").Value = 2282
$ws.Range("B9").Value = 1616
$ws.Range("C7").Value = 6365
$ws.Range("C8").Value = 6110
$ws.Range("C9").Value = 5023

# --- Column D width (best-fit like Excel produced for the new numbers) ---
$ws.Columns.Item(4).ColumnWidth = 11

# --- View state: scroll + selection ---------------------------------------
$ws.Range("I8:J8").Select() | Out-Null
